# Updates the two validation-report worksheets ("Summary" and "Details") to
# reflect the regenerated report after reverting to the old V3 SSOC assigner
# (GMI/HQA filtering removed). New rule rows (RULE 17, RULE 18, RULE 19) are
# introduced, some rule messages are reassigned to different rows, and the
# Details sheet gains corresponding per-response rows; a few existing
# member_index cells are cleared to blank.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Summary" sheet: columns A=rule, B=column, C=message, D=count
# ---------------------------------------------------------------------
$summaryWs = $wb.Worksheets.Item("Summary")

$summaryRows = @(
    @('RULE 1 - Place of Birth', 'Place of Birth', 'Others answer approved with RSPD confirmation (original word count: 4, now meets minimum requirement)', 2),
    @('RULE 18', 'Place of Birth', 'Invalid country in Others: Place of Birth', 2),
    @('RULE 7', 'Employment Status as of last week & Did you perform any freelance or assignment-based work via any of the following online platform(s) in the last 12 months?', 'Mismatch: Freelance work selected but Employment Status is not Own Account Worker.', 2),
    @('RULE 1 - Place of Birth', 'Place of Birth', 'Others answer matches predefined option: ''India''', 1),
    @('RULE 10', 'Was your main job last week a paid internship, traineeship or apprenticeship? & Type of Employment?', 'Internship/Traineeship/Apprenticeship must be Fixed-Term contract employee', 1),
    @('RULE 17', 'What is your religion?', 'Normalized to ''No religion''', 1),
    @('RULE 19', 'At any point in the last 12 months, were you self-employed? & At any point in the last 12 months, did you work on your own (i.e., without paid employees) while running your own business or trade? & Did you perform any freelance or assignment-based work via any of the following online platform(s) in the last 12 months?', 'Freelance selected but self-employed/own-account not both Yes', 1),
    @('RULE 5', 'How much interest did you receive from savings (e.g., current and saving accounts, fixed deposits) in the last 12 months?', 'Invalid interest. Must be numeric between 0 and 10 (decimals allowed).', 1)
)

$r = 2
foreach ($row in $summaryRows) {
    $summaryWs.Cells.Item($r, 1).Value = $row[0]
    $summaryWs.Cells.Item($r, 2).Value = $row[1]
    $summaryWs.Cells.Item($r, 3).Value = $row[2]
    $summaryWs.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# "Details" sheet: columns A=file, B=row, C=response_id, D=member_index,
# E=member, F=rule, G=column, H=message
# ---------------------------------------------------------------------
$detailsWs = $wb.Worksheets.Item("Details")

$detailsRows = @(
    @('CLFS_sample_input.xlsx', 2, '697c21e233bb9b02fa208b14', $null, 'handsome', 'RULE 17', 'What is your religion?', 'Normalized to ''No religion'''),
    @('CLFS_sample_input.xlsx', 3, '697c268c100972b18200a4cb', $null, 'Donald Trump', 'RULE 18', 'Place of Birth', 'Invalid country in Others: Place of Birth'),
    @('CLFS_sample_input.xlsx', 5, '697c390afe776e2d3cb543d7', $null, 'Donald Trump', 'RULE 18', 'Place of Birth', 'Invalid country in Others: Place of Birth'),
    @('CLFS_sample_input.xlsx', 2, '697c21e233bb9b02fa208b14', $null, 'handsome', 'RULE 1 - Place of Birth', 'Place of Birth', 'Others answer matches predefined option: ''India'''),
    @('CLFS_sample_input.xlsx', 3, '697c268c100972b18200a4cb', $null, 'Donald Trump', 'RULE 1 - Place of Birth', 'Place of Birth', 'Others answer approved with RSPD confirmation (original word count: 4, now meets minimum requirement)'),
    @('CLFS_sample_input.xlsx', 5, '697c390afe776e2d3cb543d7', $null, 'Donald Trump', 'RULE 1 - Place of Birth', 'Place of Birth', 'Others answer approved with RSPD confirmation (original word count: 4, now meets minimum requirement)'),
    @('CLFS_sample_input.xlsx', 2, '697c21e233bb9b02fa208b14', 1, 'handsome', 'RULE 5', 'How much interest did you receive from savings (e.g., current and saving accounts, fixed deposits) in the last 12 months?', 'Invalid interest. Must be numeric between 0 and 10 (decimals allowed).'),
    @('CLFS_sample_input.xlsx', 2, '697c21e233bb9b02fa208b14', 1, 'handsome', 'RULE 7', 'Employment Status as of last week & Did you perform any freelance or assignment-based work via any of the following online platform(s) in the last 12 months?', 'Mismatch: Freelance work selected but Employment Status is not Own Account Worker.'),
    @('CLFS_sample_input.xlsx', 2, '697c21e233bb9b02fa208b14', 1, 'handsome', 'RULE 19', 'At any point in the last 12 months, were you self-employed? & At any point in the last 12 months, did you work on your own (i.e., without paid employees) while running your own business or trade? & Did you perform any freelance or assignment-based work via any of the following online platform(s) in the last 12 months?', 'Freelance selected but self-employed/own-account not both Yes'),
    @('CLFS_sample_input.xlsx', 4, '697c2c580deae81fbb49c180', 1, 'Chen Jia Hui', 'RULE 7', 'Employment Status as of last week & Did you perform any freelance or assignment-based work via any of the following online platform(s) in the last 12 months?', 'Mismatch: Freelance work selected but Employment Status is not Own Account Worker.'),
    @('CLFS_sample_input.xlsx', 4, '697c2c580deae81fbb49c180', 1, 'Chen Jia Hui', 'RULE 10', 'Was your main job last week a paid internship, traineeship or apprenticeship? & Type of Employment?', 'Internship/Traineeship/Apprenticeship must be Fixed-Term contract employee')
)

$r = 2
foreach ($row in $detailsRows) {
    $detailsWs.Cells.Item($r, 1).Value = $row[0]
    $detailsWs.Cells.Item($r, 2).Value = $row[1]
    $detailsWs.Cells.Item($r, 3).Value = $row[2]
    if ($null -eq $row[3]) {
        $detailsWs.Cells.Item($r, 4).Value = ""
    } else {
        $detailsWs.Cells.Item($r, 4).Value = $row[3]
    }
    $detailsWs.Cells.Item($r, 5).Value = $row[4]
    $detailsWs.Cells.Item($r, 6).Value = $row[5]
    $detailsWs.Cells.Item($r, 7).Value = $row[6]
    $detailsWs.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}
